$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.832.39"
$ws.Range("E2").Value = "  +5.67%  "
$ws.Range("D3").Value = "'3.111.61"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("D5").Value = "'585.65"
$ws.Range("E5").Value = "  +3.59%  "
$ws.Range("D6").Value = "'143.58"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'3.100.99"
$ws.Range("E8").Value = "  +3.33%  "
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("E10").Value = "  +8.79%  "
$ws.Range("D11").Value = "'5.76"
$ws.Range("E11").Value = "  +10.05%  "
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("E13").Value = "  +4.67%  "
$ws.Range("D14").Value = "'35.54"
$ws.Range("E14").Value = "  +4.82%  "
$ws.Range("D15").Value = "'0.123"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "'3.625.79"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'62.781.45"
$ws.Range("E18").Value = "  +5.65%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "'3.102.83"
$ws.Range("E19").Value = "  +3.15%  "
$ws.Range("D20").Value = "'452.21"
$ws.Range("E20").Value = "  +4.22%  "
$ws.Range("D21").Value = "'14.07"
$ws.Range("E21").Value = "  +2.95%  "
$ws.Range("D22").Value = "'0.734"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").Value = "'7.53"
$ws.Range("E23").Value = "  +5.12%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "'81.68"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +5.30%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +4.11%  "
$ws.Range("D31").Value = "'6.85"
$ws.Range("E31").Value = "  +11.33%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'27.12"
$ws.Range("E32").Value = "  +4.82%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.111"
$ws.Range("E33").Value = "  +11.03%  "
$ws.Range("D34").Value = "'1.04"
$ws.Range("E34").Value = "  +3.81%  "
$ws.Range("E35").Value = "  +4.92%  "
$ws.Range("D36").Value = "'6.07"
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").Value = "'2.27"
$ws.Range("E37").Value = "  +6.76%  "
$ws.Range("D38").Value = "'50.84"
$ws.Range("E38").Value = "  +3.74%  "
$ws.Range("D39").Value = "'3.02"
$ws.Range("E39").Value = "  +8.85%  "
$ws.Range("D40").Value = "'8.79"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("D41").Value = "'426.12"
$ws.Range("E41").Value = "  +3.84%  "
$ws.Range("D42").Value = "'2.953.51"
$ws.Range("E42").Value = "  +6.45%  "
$ws.Range("E43").Value = "  +5.15%  "
$ws.Range("D44").Value = "'0.278"
$ws.Range("E44").Value = "  +9.43%  "
$ws.Range("E45").Value = "  +2.80%  "
$ws.Range("E46").Value = "  +6.42%  "
$ws.Range("D47").Value = "'125.79"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("D49").Value = "'35.13"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("D51").Value = "'24.80"
$ws.Range("E51").Value = "  +4.91%  "
